# Scheduled runner refresh: update market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the per-class Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H19").Value = 716042.9
$ws_ALC.Range("I19").Value = 1251453.5
$ws_ALC.Range("K19").Value = 1251453.5
$ws_ALC.Range("M19").Value = -1251278.5
$ws_ALC.Range("H33").Value = 120.666664
$ws_ALC.Range("H62").Value = 10394.667
$ws_ALC.Range("I62").Value = 8097
$ws_ALC.Range("K62").Value = 8097
$ws_ALC.Range("M62").Value = -7473
$ws_ALC.Range("H65").Value = 10394.667
$ws_ALC.Range("I65").Value = 8097
$ws_ALC.Range("K65").Value = 40485
$ws_ALC.Range("M65").Value = -37365
$ws_ALC.Range("H70").Value = 13420.375
$ws_ALC.Range("I70").Value = 6149.8335
$ws_ALC.Range("K70").Value = 18449.5005
$ws_ALC.Range("M70").Value = -18179.5005
$ws_ALC.Range("H73").Value = 13420.375
$ws_ALC.Range("I73").Value = 6149.8335
$ws_ALC.Range("K73").Value = 18449.5005
$ws_ALC.Range("M73").Value = -17513.5005
$ws_ALC.Range("H97").Value = 2033
$ws_ALC.Range("J97").Value = 2033
$ws_ALC.Range("L97").Value = 6099
$ws_ALC.Range("N97").Value = -7091
$ws_ALC.Range("H132").Value = 885.1070999999999
$ws_ALC.Range("I132").Value = 876.2692
$ws_ALC.Range("K132").Value = 2628.8076
$ws_ALC.Range("M132").Value = -98.80760000000009
$ws_ALC.Range("H138").Value = 2831.5557
$ws_ALC.Range("J138").Value = 1929.0385
$ws_ALC.Range("L138").Value = 5787.1155
$ws_ALC.Range("N138").Value = -16067.1155
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 3330.3135
$ws_ARM.Range("I32").Value = 2785.7925
$ws_ARM.Range("K32").Value = 2785.7925
$ws_ARM.Range("M32").Value = -2498.7925
$ws_ARM.Range("H45").Value = 2265.2812
$ws_ARM.Range("I45").Value = 1570.6428
$ws_ARM.Range("J45").Value = 2805.5557
$ws_ARM.Range("K45").Value = 1570.6428
$ws_ARM.Range("L45").Value = 2805.5557
$ws_ARM.Range("M45").Value = -1193.6428
$ws_ARM.Range("N45").Value = -3559.5557
$ws_ARM.Range("H74").Value = 775.0476
$ws_ARM.Range("I74").Value = 664
$ws_ARM.Range("J74").Value = 1247
$ws_ARM.Range("K74").Value = 664
$ws_ARM.Range("L74").Value = 1247
$ws_ARM.Range("M74").Value = 210
$ws_ARM.Range("N74").Value = -2995
$ws_ARM.Range("H77").Value = 775.0476
$ws_ARM.Range("I77").Value = 664
$ws_ARM.Range("J77").Value = 1247
$ws_ARM.Range("K77").Value = 3320
$ws_ARM.Range("L77").Value = 6235
$ws_ARM.Range("M77").Value = 1048
$ws_ARM.Range("N77").Value = -14971
$ws_ARM.Range("H132").Value = 4265.3335
$ws_ARM.Range("I132").Value = 3898.5
$ws_ARM.Range("K132").Value = 11695.5
$ws_ARM.Range("M132").Value = -9165.5
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H105").Value = 2572.55
$ws_BSM.Range("I105").Value = 2586.8948
$ws_BSM.Range("J105").Value = 2300
$ws_BSM.Range("K105").Value = 2586.8948
$ws_BSM.Range("L105").Value = 2300
$ws_BSM.Range("M105").Value = -839.8948
$ws_BSM.Range("N105").Value = -5794
$ws_BSM.Range("H134").Value = 5613.3335
$ws_BSM.Range("J134").Value = 3957
$ws_BSM.Range("L134").Value = 11871
$ws_BSM.Range("N134").Value = -16941
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1493.0741
$ws_CRP.Range("I31").Value = 958.5714
$ws_CRP.Range("J31").Value = 2068.6924
$ws_CRP.Range("K31").Value = 958.5714
$ws_CRP.Range("L31").Value = 2068.6924
$ws_CRP.Range("M31").Value = -663.5714
$ws_CRP.Range("N31").Value = -2658.6924
$ws_CRP.Range("H34").Value = 1493.0741
$ws_CRP.Range("I34").Value = 958.5714
$ws_CRP.Range("J34").Value = 2068.6924
$ws_CRP.Range("K34").Value = 958.5714
$ws_CRP.Range("L34").Value = 2068.6924
$ws_CRP.Range("M34").Value = -756.5714
$ws_CRP.Range("N34").Value = -2472.6924
$ws_CRP.Range("H43").Value = 14883.667
$ws_CRP.Range("J43").Value = 14883.667
$ws_CRP.Range("L43").Value = 14883.667
$ws_CRP.Range("N43").Value = -15251.667
$ws_CRP.Range("H62").Value = 9993.333000000001
$ws_CRP.Range("I62").Value = 9993.333000000001
$ws_CRP.Range("K62").Value = 9993.333000000001
$ws_CRP.Range("M62").Value = -9369.333000000001
$ws_CRP.Range("H65").Value = 9993.333000000001
$ws_CRP.Range("I65").Value = 9993.333000000001
$ws_CRP.Range("K65").Value = 49966.665
$ws_CRP.Range("M65").Value = -46846.665
$ws_CRP.Range("H101").Value = 14883.667
$ws_CRP.Range("J101").Value = 14883.667
$ws_CRP.Range("L101").Value = 14883.667
$ws_CRP.Range("N101").Value = -21373.667
$ws_CRP.Range("H105").Value = 1701.125
$ws_CRP.Range("I105").Value = 1701.125
$ws_CRP.Range("K105").Value = 1701.125
$ws_CRP.Range("M105").Value = 45.875
$ws_CRP.Range("H134").Value = 1901.7646
$ws_CRP.Range("I134").Value = 1240
$ws_CRP.Range("K134").Value = 3720
$ws_CRP.Range("M134").Value = -1185
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H41").Value = 6777.6
$ws_CUL.Range("I41").Value = 3000
$ws_CUL.Range("K41").Value = 9000
$ws_CUL.Range("M41").Value = -8662
$ws_CUL.Range("H131").Value = 20863884
$ws_CUL.Range("I131").Value = 50000404
$ws_CUL.Range("J131").Value = 52084.93
$ws_CUL.Range("K131").Value = 150001212
$ws_CUL.Range("L131").Value = 156254.79
$ws_CUL.Range("M131").Value = -149996172
$ws_CUL.Range("N131").Value = -166334.79
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 3030.4
$ws_GSM.Range("I102").Value = 3288.75
$ws_GSM.Range("J102").Value = 1997
$ws_GSM.Range("K102").Value = 3288.75
$ws_GSM.Range("L102").Value = 1997
$ws_GSM.Range("M102").Value = -1666.75
$ws_GSM.Range("N102").Value = -5241
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H55").Value = 130.3077
$ws_LTW.Range("I55").Value = 144.45454
$ws_LTW.Range("J55").Value = 52.5
$ws_LTW.Range("K55").Value = 144.45454
$ws_LTW.Range("L55").Value = 52.5
$ws_LTW.Range("M55").Value = 28.54545999999999
$ws_LTW.Range("N55").Value = -398.5
$ws_LTW.Range("H81").Value = 49249.5
$ws_LTW.Range("J81").Value = 49249.5
$ws_LTW.Range("L81").Value = 49249.5
$ws_LTW.Range("N81").Value = -51245.5
$ws_LTW.Range("H84").Value = 49249.5
$ws_LTW.Range("J84").Value = 49249.5
$ws_LTW.Range("L84").Value = 147748.5
$ws_LTW.Range("N84").Value = -157732.5
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H80").Value = 50000
$ws_WVR.Range("J80").Value = 50000
$ws_WVR.Range("L80").Value = 50000
$ws_WVR.Range("N80").Value = -51996
$ws_WVR.Range("H83").Value = 50000
$ws_WVR.Range("J83").Value = 50000
$ws_WVR.Range("L83").Value = 150000
$ws_WVR.Range("N83").Value = -159984
$ws_WVR.Range("H122").Value = 145889.08
$ws_WVR.Range("I122").Value = 157838.25
$ws_WVR.Range("J122").Value = 2499
$ws_WVR.Range("K122").Value = 473514.75
$ws_WVR.Range("L122").Value = 7497
$ws_WVR.Range("M122").Value = -471064.75
$ws_WVR.Range("N122").Value = -12397
$ws_WVR.Range("H132").Value = 1347.3529
$ws_WVR.Range("I132").Value = 565.5
$ws_WVR.Range("J132").Value = 4996
$ws_WVR.Range("K132").Value = 1696.5
$ws_WVR.Range("L132").Value = 14988
$ws_WVR.Range("M132").Value = 833.5
$ws_WVR.Range("N132").Value = -20048
$ws_WVR.Range("H136").Value = 42739732
$ws_WVR.Range("I136").Value = 92598350
$ws_WVR.Range("J136").Value = 3778.4285
$ws_WVR.Range("K136").Value = 277795050
$ws_WVR.Range("L136").Value = 11335.2855
$ws_WVR.Range("M136").Value = -277792500
$ws_WVR.Range("N136").Value = -16435.2855
